$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the CLASSIFICATION header (C1) text to INFORMATIONDOMAIN.
$ws.Range("C1").Value = "INFORMATIONDOMAIN"

# Swap the styles between B1/C1 and update B2 to match.
$orgStyle = $ws.Range("B1").Style
$classStyle = $ws.Range("C1").Style
$ws.Range("B1").Style = $classStyle
$ws.Range("C1").Style = $orgStyle
$ws.Range("B2").Style = $classStyle
